$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 302.46155
$ws.Range("J33").Value = 57
$ws.Range("L33").Value = 57
$ws.Range("N33").Value = -515

$ws.Range("H64").Value = 250004510
$ws.Range("J64").Value = 500004300
$ws.Range("L64").Value = 500004300
$ws.Range("N64").Value = -500004796

$ws.Range("H67").Value = 250004510
$ws.Range("J67").Value = 500004300
$ws.Range("L67").Value = 500004300
$ws.Range("N67").Value = -500006016

$ws.Range("H106").Value = 1202.25

$ws.Range("H112").Value = 1629.2941
$ws.Range("I112").Value = 461.66666
$ws.Range("J112").Value = 1879.5
$ws.Range("K112").Value = 1384.99998
$ws.Range("L112").Value = 5638.5
$ws.Range("M112").Value = -276.9999800000001
$ws.Range("N112").Value = -7854.5

$ws.Range("H125").Value = 957.6667
$ws.Range("I125").Value = 935.2857
$ws.Range("K125").Value = 8417.5713
$ws.Range("M125").Value = -5957.5713

$ws.Range("H137").Value = 1285772.8
$ws.Range("I137").Value = 1926243.5
$ws.Range("J137").Value = 4831.3076
$ws.Range("K137").Value = 5778730.5
$ws.Range("L137").Value = 14493.9228
$ws.Range("M137").Value = -5776180.5
$ws.Range("N137").Value = -19593.9228

$ws.Range("H138").Value = 2273.3333
$ws.Range("I138").Value = 1102
$ws.Range("K138").Value = 3306
$ws.Range("M138").Value = 1834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1568837.1
$ws.Range("I32").Value = 713387.7
$ws.Range("K32").Value = 713387.7
$ws.Range("M32").Value = -713100.7

$ws.Range("H63").Value = 1947
$ws.Range("I63").Value = 1894
$ws.Range("K63").Value = 1894
$ws.Range("M63").Value = -1208

$ws.Range("H66").Value = 1947
$ws.Range("I66").Value = 1894
$ws.Range("K66").Value = 9470
$ws.Range("M66").Value = -6038

$ws.Range("H122").Value = 2693.6843
$ws.Range("I122").Value = 2579.4243
$ws.Range("J122").Value = 3447.8
$ws.Range("K122").Value = 7738.2729
$ws.Range("L122").Value = 10343.4
$ws.Range("M122").Value = -5288.2729
$ws.Range("N122").Value = -15243.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 993.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 993.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 993.5
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -2989.5

$ws.Range("H83").Value = 993.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 993.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 4967.5
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -14951.5

$ws.Range("H86").Value = 3849.4666
$ws.Range("I86").Value = 3660.1428
$ws.Range("K86").Value = 3660.1428
$ws.Range("M86").Value = -2537.1428

$ws.Range("H89").Value = 3849.4666
$ws.Range("I89").Value = 3660.1428
$ws.Range("K89").Value = 18300.714
$ws.Range("M89").Value = -12684.714

$ws.Range("H94").Value = 133340920
$ws.Range("I94").Value = 222234020
$ws.Range("K94").Value = 222234020
$ws.Range("M94").Value = -222233569

$ws.Range("H105").Value = 13001638
$ws.Range("J105").Value = 25001720
$ws.Range("L105").Value = 25001720
$ws.Range("N105").Value = -25005214

$ws.Range("H107").Value = 3498239.5
$ws.Range("I107").Value = 4275220
$ws.Range("K107").Value = 4275220
$ws.Range("M107").Value = -4273300

$ws.Range("H134").Value = 1193.1714
$ws.Range("I134").Value = 847.1613
$ws.Range("K134").Value = 2541.4839
$ws.Range("M134").Value = -6.483900000000176

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 3500
$ws.Range("J8").Value = 3500
$ws.Range("L8").Value = 3500
$ws.Range("N8").Value = -3780

$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5348

$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650

$ws.Range("H31").Value = 1788087
$ws.Range("I31").Value = 1505.7916
$ws.Range("J31").Value = 2720216.2
$ws.Range("K31").Value = 1505.7916
$ws.Range("L31").Value = 2720216.2
$ws.Range("M31").Value = -1210.7916
$ws.Range("N31").Value = -2720806.2

$ws.Range("H34").Value = 1788087
$ws.Range("I34").Value = 1505.7916
$ws.Range("J34").Value = 2720216.2
$ws.Range("K34").Value = 1505.7916
$ws.Range("L34").Value = 2720216.2
$ws.Range("M34").Value = -1303.7916
$ws.Range("N34").Value = -2720620.2

$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H134").Value = 3756.2917
$ws.Range("I134").Value = 3797.762
$ws.Range("K134").Value = 11393.286
$ws.Range("M134").Value = -8858.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1490
$ws.Range("I51").Value = 1490
$ws.Range("K51").Value = 4470
$ws.Range("M51").Value = -4010

$ws.Range("H129").Value = 1666.1428
$ws.Range("I129").Value = 1336
$ws.Range("J129").Value = 2491.5
$ws.Range("K129").Value = 4008
$ws.Range("L129").Value = 7474.5
$ws.Range("M129").Value = 992
$ws.Range("N129").Value = -17474.5

$ws.Range("H131").Value = 12501991
$ws.Range("J131").Value = 1823.7693
$ws.Range("L131").Value = 5471.3079
$ws.Range("N131").Value = -15551.3079

$ws.Range("H140").Value = 2168.1428
$ws.Range("I140").Value = 2168.1428
$ws.Range("K140").Value = 6504.428400000001
$ws.Range("M140").Value = -1324.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 23106.834
$ws.Range("J98").Value = 23106.834
$ws.Range("L98").Value = 23106.834
$ws.Range("N98").Value = -29096.834

$ws.Range("H107").Value = 9515.666999999999
$ws.Range("J107").Value = 16736.5
$ws.Range("L107").Value = 16736.5
$ws.Range("N107").Value = -20576.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 76531940
$ws.Range("I22").Value = 7938129
$ws.Range("K22").Value = 7938129
$ws.Range("M22").Value = -7937834

$ws.Range("H27").Value = 76531940
$ws.Range("I27").Value = 7938129
$ws.Range("K27").Value = 7938129
$ws.Range("M27").Value = -7938022

$ws.Range("H122").Value = 7234.5
$ws.Range("I122").Value = 5900.154
$ws.Range("K122").Value = 17700.462
$ws.Range("M122").Value = -15250.462

$ws.Range("H132").Value = 4351.5293
$ws.Range("I132").Value = 4279.636
$ws.Range("K132").Value = 12838.908
$ws.Range("M132").Value = -10308.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws.Range("H105").Value = 17615
$ws.Range("J105").Value = 17615
$ws.Range("L105").Value = 17615
$ws.Range("N105").Value = -24603

$ws.Range("H126").Value = 2241.3333
$ws.Range("I126").Value = 2241.3333
$ws.Range("K126").Value = 6723.999899999999
$ws.Range("M126").Value = -4253.999899999999
